$wb = $excel.ActiveWorkbook

# New sheet "PLan5" is a copy of Plan4 (same "ano/custos/receitas" layout),
# inserted right after Plan4 so it becomes sheetId=5 / rId5, the last tab,
# and (per Excel's normal behaviour) the new active sheet.
$plan4 = $wb.Worksheets.Item(4)
$plan4.Copy($null, $plan4) | Out-Null
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "PLan5"

# The copy starts out identical to Plan4 (headers + the 0..7 "ano" column);
# fill in the extra "custos"/"receitas" data points that differ on PLan5.
$ws.Range("B2").Value = 2500
$ws.Range("B9").Value = 1000
$ws.Range("C9").Value = 10000

# Selection on the new (now active) sheet.
$ws.Range("B10").Select() | Out-Null
